# Finished the presentation slides.
# Adds two new rows (exam dates) to the "2017-2018" sheet, below the
# existing weekly-plan table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "2017-2018"

# --- Row 29: Exame 1 --------------------------------------------------
$ws.Range("A29").Value = 43467
$ws.Range("A29").NumberFormat = "DD/MM/YYYY"

$ws.Range("B29").Value = "Exame 1"
$ws.Range("B29").HorizontalAlignment = -4108   # xlCenter

$ws.Range("C29").Value = 0.479166666666667
$ws.Range("C29").NumberFormat = "HH:MM:SS"
$ws.Range("C29").HorizontalAlignment = -4108   # xlCenter

# --- Row 30: Exame 2 --------------------------------------------------
$ws.Range("A30").Value = 43501
$ws.Range("A30").NumberFormat = "DD/MM/YYYY"

$ws.Range("B30").Value = "Exame 2"
$ws.Range("B30").HorizontalAlignment = -4108   # xlCenter

$ws.Range("C30").Value = 0.479166666666667
$ws.Range("C30").NumberFormat = "HH:MM:SS"
$ws.Range("C30").HorizontalAlignment = -4108   # xlCenter

# Move / record the selection where the user left off, like the source
# workbook shows (one row below the new last row).
$ws.Range("B31").Select()
